$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Row 4 data updates (valorPago, numeroPrestamo, numeroCuenta, tipoCuenta, tipoPrestamo)
$ws.Range("R4").Value = "1000"
$ws.Range("P4").Value = "29281005233"
$ws.Range("T4").Value = "406-107870-00"
$ws.Range("O4").Value = "CREDIAGIL"
$ws.Range("S4").Value = "Corriente"

# Sheet view: scroll position and selection changed
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("P4").Select()
